$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 2: Ecriture de la documentation - duration 3 -> 8
$ws.Range("D2").Value = 8

# Row 5: Affichage des produits, gestion par administrateur - start date 42122 -> 42124, duration blank -> 1
$ws.Range("C5").Value = 42124
$ws.Range("D5").Value = 1

# Row 7: Affichage des produits par mot clés - start date 42122 -> 42124, duration blank -> 5
$ws.Range("C7").Value = 42124
$ws.Range("D7").Value = 5

# Row 10: Affichage détaillé d'un produit - start date 42122 -> 42124, duration blank -> 6
$ws.Range("C10").Value = 42124
$ws.Range("D10").Value = 6

# Update view selection
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("I11").Select()
